# Update "想去人数" (want-to-go count) and occasionally "最低票价" (min price)
# figures across the four sheets to match the newly scraped snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1165
$ws.Range("F3").Value = 1071
$ws.Range("F4").Value = 1867
$ws.Range("F5").Value = 588
$ws.Range("F6").Value = 1229
$ws.Range("F7").Value = 63
$ws.Range("F8").Value = 22
$ws.Range("F9").Value = 131
$ws.Range("F10").Value = 321
$ws.Range("F11").Value = 98
$ws.Range("F12").Value = 95
$ws.Range("F13").Value = 764
$ws.Range("F14").Value = 209
$ws.Range("F15").Value = 118
$ws.Range("G15").Value = 60
$ws.Range("F18").Value = 336
$ws.Range("F19").Value = 188
$ws.Range("F20").Value = 688
$ws.Range("F22").Value = 658
$ws.Range("F23").Value = 177
$ws.Range("F24").Value = 41
$ws.Range("F25").Value = 890
$ws.Range("F26").Value = 334
$ws.Range("F27").Value = 177
$ws.Range("F29").Value = 295
$ws.Range("F30").Value = 14
$ws.Range("F32").Value = 416

# --- Sheet "演出" (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 327

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 320

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 320
$ws.Range("F3").Value = 1165
$ws.Range("F4").Value = 1071
$ws.Range("F5").Value = 1867
$ws.Range("F6").Value = 588
$ws.Range("F7").Value = 1229
$ws.Range("F8").Value = 63
$ws.Range("F10").Value = 22
$ws.Range("F11").Value = 131
$ws.Range("F12").Value = 321
$ws.Range("F13").Value = 98
$ws.Range("F14").Value = 95
$ws.Range("F15").Value = 764
$ws.Range("F16").Value = 209
$ws.Range("F17").Value = 118
$ws.Range("G17").Value = 60
$ws.Range("F20").Value = 327
$ws.Range("F23").Value = 336
$ws.Range("F27").Value = 188
$ws.Range("F28").Value = 688
$ws.Range("F30").Value = 658
$ws.Range("F31").Value = 177
$ws.Range("F32").Value = 41
$ws.Range("F33").Value = 890
$ws.Range("F34").Value = 334
$ws.Range("F37").Value = 177
$ws.Range("F39").Value = 295
$ws.Range("F43").Value = 14
$ws.Range("F46").Value = 416
